# Applies the "added date modified and additional version information" edit
# to "Ra for Android.docx".
#
# Strategy: use Range.InsertXML with small WordprocessingML fragments to get
# exact control over run/paragraph boundaries (run-splits, field codes,
# bookmarks, lastRenderedPageBreak, etc.), and do the structural edits from
# the bottom of the document upward so earlier paragraph indices stay valid.

$d = $word.ActiveDocument
$wNs = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) "Version 2.0" section -> expanded version history (0.2, 0.3, 1.0)
#    Replaces paragraphs 23 ("Version 2.0") and 24 (the body paragraph
#    that follows it, which also holds the _GoBack bookmark).
# ---------------------------------------------------------------------
$pHeading20 = $d.Paragraphs(23)
$pBody20 = $d.Paragraphs(24)
$rng = $d.Range($pHeading20.Range.Start, $pBody20.Range.End)

$xml = "<w:p$wNs><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr>" `
     + "<w:r><w:t xml:space=`"preserve`">Version </w:t></w:r>" `
     + "<w:r><w:t>0.2</w:t></w:r></w:p>"
$xml += "<w:p$wNs>" `
     + "<w:r><w:t xml:space=`"preserve`">The second version </w:t></w:r>" `
     + "<w:r><w:t>will/</w:t></w:r>" `
     + "<w:r><w:t>was my project for Winter quarter 2017. The object is to greatly improve the visual interface using graphics and other UI functionality that will be taught in the class for the winter quarter.</w:t></w:r>" `
     + "</w:p>"
$xml += "<w:p$wNs><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr>" `
     + "<w:r><w:t>Version 0.3</w:t></w:r></w:p>"
$xml += "<w:p$wNs>" `
     + "<w:r><w:t xml:space=`"preserve`">The third version will/did use </w:t></w:r>" `
     + "</w:p>"
$xml += "<w:p$wNs><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr>" `
     + "<w:r><w:lastRenderedPageBreak/><w:t>Version 1.0</w:t></w:r></w:p>"
$xml += "<w:p$wNs>" `
     + "<w:r><w:t>The first major version will be to make a commercial quality app, to be placed online for download. May charge for app, but will at least be</w:t></w:r>" `
     + "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" `
     + "<w:r><w:t xml:space=`"preserve`"> a demonstration of my work to potential employers.</w:t></w:r>" `
     + "</w:p>"

$rng.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------
# 2) First "Version 1.0" heading -> "Version " + "0.1" (two runs)
# ---------------------------------------------------------------------
$pHeading10 = $d.Paragraphs(8)
$xml = "<w:p$wNs><w:pPr><w:pStyle w:val=`"Heading2`"/></w:pPr>" `
     + "<w:r><w:t xml:space=`"preserve`">Version </w:t></w:r>" `
     + "<w:r><w:t>0.1</w:t></w:r></w:p>"
$pHeading10.Range.InsertXML($xml) | Out-Null

# ---------------------------------------------------------------------
# 3) "By Benjamin Sklar" byline paragraph -> NoSpacing style, single run,
#    plus a new "Last Modified: <DATE field>" paragraph right after it.
# ---------------------------------------------------------------------
$pByline = $d.Paragraphs(2)

# Mint the built-in "No Spacing" paragraph style (matches Word's own
# NoSpacing/"No Spacing" style: uiPriority 1, qFormat, spacing after=0,
# single line spacing) by assigning it as a paragraph style first.
$pByline.Style = "No Spacing"
$noSpacing = $d.Styles("No Spacing")
$noSpacing.Priority = 1
$noSpacingPf = $noSpacing.ParagraphFormat
$noSpacingPf.SpaceAfter = 0
$noSpacingPf.LineSpacingRule = 0

$xml = "<w:p$wNs><w:pPr><w:pStyle w:val=`"NoSpacing`"/></w:pPr>" `
     + "<w:r><w:t>By Benjamin Sklar</w:t></w:r></w:p>"
$xml += "<w:p$wNs>" `
     + "<w:r><w:t xml:space=`"preserve`">Last Modified: </w:t></w:r>" `
     + "<w:r><w:fldChar w:fldCharType=`"begin`"/></w:r>" `
     + "<w:r><w:instrText xml:space=`"preserve`"> DATE \@ `"yyyy-MM-dd`" </w:instrText></w:r>" `
     + "<w:r><w:fldChar w:fldCharType=`"separate`"/></w:r>" `
     + "<w:r><w:rPr><w:noProof/></w:rPr><w:t>2017-01-04</w:t></w:r>" `
     + "<w:r><w:fldChar w:fldCharType=`"end`"/></w:r>" `
     + "</w:p>"
$pByline.Range.InsertXML($xml) | Out-Null

Write-Output "edit complete"
